$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row (for "Brasil" origin Mango, bandeja 160 units,
# date serial 44777) is inserted at row 140, pushing the existing rows
# 140-222 down to 141-223 (dimension grows from A1:T222 to A1:T223).
$ws.Rows.Item(140).Insert()

$ws.Cells.Item(140, 1).Value  = 4
$ws.Cells.Item(140, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(140, 3).Value  = "Los Lagos"
$ws.Cells.Item(140, 4).Value  = 44777
$ws.Cells.Item(140, 5).Value  = 10
$ws.Cells.Item(140, 6).Value  = "Fruta"
$ws.Cells.Item(140, 7).Value  = 100108
$ws.Cells.Item(140, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(140, 9).Value  = 100108002
$ws.Cells.Item(140, 10).Value = "Mango"
$ws.Cells.Item(140, 11).Value = "Sin especificar"
$ws.Cells.Item(140, 12).Value = "Primera"
$ws.Cells.Item(140, 13).Value = 160
$ws.Cells.Item(140, 14).Value = 13000
$ws.Cells.Item(140, 15).Value = 14000
$ws.Cells.Item(140, 16).Value = 13500
$ws.Cells.Item(140, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(140, 18).Value = "Brasil"
$ws.Cells.Item(140, 19).Value = 3375
$ws.Cells.Item(140, 20).Value = 4
